$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values regenerated after filtering save games (rows correspond to dates in column A)
# Columns: B=TB, C=d2S, D=K, E=IP, F=Win (unchanged), G=sum

$data = @{
    2 = @{ B = 0.04763786555579896; C = 0.04240448674262143; D = 0.8054896365839992;  E = 0.496779210170732; G = 1.392311199053152 }
    3 = @{ B = 3.230985683306322;   C = 1.667794583268128;  D = 0.8054896365839992;  E = 0.496779210170732; G = 6.201049113329182 }
    4 = @{ B = 0.127881588408715;   C = 0.3127903958511391; D = 3.900430680208489;   E = 0.496779210170732; G = 4.837881874639075 }
    5 = @{ B = 1.459612070389937;   C = 1.667794583268128;  D = 0.1575252929769615;  E = 0.496779210170732; G = 3.781711156805759 }
    6 = @{ B = 3.230985683306322;   C = 1.667794583268128;  D = 0.8054896365839992;  E = 0.496779210170732; G = 6.201049113329182 }
    7 = @{ B = 3.230985683306322;   C = 1.667794583268128;  D = 3.900430680208489;   E = 0.496779210170732; G = 9.295990156953671 }
    8 = @{ B = 3.230985683306322;   C = 10.29869402782916;  D = 0.1575252929769615;  E = 8.660232485948974; G = 22.34743749006142 }
    9 = @{ B = 0.6753301551942219;  C = 1.667794583268128;  D = 0.8054896365839992;  E = 0.496779210170732; G = 3.645393585217082 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
